$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I1 and J1, copying the formatting from the existing
# header cell H1 (bold font, thin border, centered alignment) so the
# new header cells share style index 1 just like the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the I0 / IF values for rows 2-73.
$data = @(
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(8,8),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(8,8),
    @(9,9),
    @(8,8),
    @(9,9),
    @(9,9),
    @(9,9),
    @(8,8),
    @(8,8),
    @(8,9),
    @(8,8),
    @(9,9),
    @(9,9),
    @(9,9),
    @(7,8),
    @(9,9),
    @(8,9),
    @(9,9),
    @(10,10),
    @(8,8),
    @(9,9),
    @(9,9),
    @(7,7),
    @(8,8),
    @(7,8),
    @(8,8),
    @(9,9),
    @(7,7),
    @(8,8),
    @(9,9),
    @(9,9),
    @(9,9),
    @(7,8),
    @(9,9),
    @(8,8),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(6,6),
    @(5,5),
    @(7,7),
    @(5,5),
    @(6,6),
    @(3,3)
)

$r = 2
foreach ($pair in $data) {
    $ws.Cells.Item($r, 9).Value = $pair[0]
    $ws.Cells.Item($r, 10).Value = $pair[1]
    $r = $r + 1
}
